$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J2").Value = 6
$ws.Range("K2").Value = "6/52"
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = "1/5"
